# Apply edits described by the diff: add two new OCXO rows (13 and 14)
# on Sheet1, and move the active cell selection to J15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlCenter = -4108

# ---------------------------------------------------------------
# Row 13: Abracon Corp / AOCJY-16.384MHZ OCXO
# (order of writes chosen to reproduce the original shared-string
#  insertion order: Abracon Corp, 0-50, 5/3.3, 25.4x22.1,
#  AOCJY-16.384MHZ)
# ---------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "OCXO"
$ws.Cells.Item(13, 2).Value = "Abracon Corp"

$ws.Cells.Item(13, 5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 5).Value = 5

$ws.Cells.Item(13, 6).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 6).Value = 20

$ws.Cells.Item(13, 7).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 7).Value = 1

$ws.Cells.Item(13, 8).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 8).Value = 0.5

$ws.Cells.Item(13, 9).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 9).Value = 10

$ws.Cells.Item(13, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 10).Value = "0-50"

$ws.Cells.Item(13, 12).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 12).Value = "5/3.3"

$ws.Cells.Item(13, 15).Value = "25.4x22.1"

$ws.Cells.Item(13, 3).Value = "AOCJY-16.384MHZ"

$ws.Cells.Item(13, 11).NumberFormat = '_("$"* #,##0.00_);_("$"* (#,##0.00);_("$"* "-"??_);_(@_)'
$ws.Cells.Item(13, 11).Value = 160.8

$ws.Cells.Item(13, 13).NumberFormat = "0"
$ws.Cells.Item(13, 13).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 13).Formula = "=N13/3.3*1000"

$ws.Cells.Item(13, 14).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 14).Value = 3.6

$ws.Cells.Item(13, 16).HorizontalAlignment = $xlCenter
$ws.Cells.Item(13, 16).Value = "LVCMOS"

# ---------------------------------------------------------------
# Row 14: Connor-Winfield / OX914xS3 OCXO
# (order of writes: OX914xS3, then -55-85)
# ---------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = "OCXO"
$ws.Cells.Item(14, 2).Value = "Connor-Winfield"
$ws.Cells.Item(14, 3).Value = "OX914xS3"

$ws.Cells.Item(14, 4).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 4).Value = 1

$ws.Cells.Item(14, 5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 5).Value = 140

$ws.Cells.Item(14, 6).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 6).Value = 20

$ws.Cells.Item(14, 7).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 7).Value = 40

$ws.Cells.Item(14, 8).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 8).Value = 4.6

$ws.Cells.Item(14, 9).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 9).Value = 20

$ws.Cells.Item(14, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 10).Value = "'-55-85"

$ws.Cells.Item(14, 11).NumberFormat = '_("$"* #,##0.00_);_("$"* (#,##0.00);_("$"* "-"??_);_(@_)'
$ws.Cells.Item(14, 11).Value = 42

$ws.Cells.Item(14, 12).Value = 3.3

$ws.Cells.Item(14, 14).Value = "1.3-3.0"

$ws.Cells.Item(14, 15).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 15).Value = "9x14"

$ws.Cells.Item(14, 16).HorizontalAlignment = $xlCenter
$ws.Cells.Item(14, 16).Value = "LVCMOS"

# --- Selection ---
$ws.Range("J15").Select()
